$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: Client Id, Candidate ID, User Name, Exam Password
$ws.Range("A2").Value = "test985"
$ws.Range("B2").Value = 23071334
$ws.Range("C2").Value = "narendra6643"
$ws.Range("D2").Value = "s%7#5hDB"
